# Weekly update: insert a new daily price record as row 10, pushing the
# remaining records (old rows 10-81) down by one row. The record that used
# to occupy the last row (81) ends up at the new last row (82) automatically
# as part of the shift.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 10.
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the new record.
$ws.Range("A10").Value = 11
$ws.Range("B10").Value = "Vega Monumental Concepción"
$ws.Range("C10").Value = "Bíobío"
$ws.Range("D10").Value = 44503
$ws.Range("E10").Value = 8
$ws.Range("F10").Value = 100112043
$ws.Range("G10").Value = "Pepino ensalada"
$ws.Range("H10").Value = "Sin especificar"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 250
$ws.Range("K10").Value = 7500
$ws.Range("L10").Value = 8000
$ws.Range("M10").Value = 7700
$ws.Range("N10").Value = "`$/caja 60 unidades"
$ws.Range("O10").Value = "Región de Arica y Parinacota"
$ws.Range("P10").Value = 128
$ws.Range("Q10").Value = 60
$ws.Range("R10").Value = "Hortaliza"
